$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Area" column (G) and its running total "Atotal" (H)
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G18").Formula = "=(D4-D3)*B4/100"

$ws.Range("H2").Formula = "=SUM(G2:G18)"

# Small summary block to the right (J:K) mirroring the totals
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Match the saved selection/view state
$ws.Range("J2:K2").Select()
